$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 22.53092045678167
$ws.Range("C2").Value = 8.721693494346196
$ws.Range("D2").Value = 12.63078510117451
$ws.Range("E2").Value = 12.21738794893011
$ws.Range("G2").Value = 58.8274148164793
$ws.Range("H2").Value = 21.69643440232802
$ws.Range("I2").Value = 33.05970174887393
$ws.Range("J2").Value = 7.685596870348577
$ws.Range("L2").Value = 13.24907472036026
$ws.Range("M2").Value = 20.38037353477413
$ws.Range("N2").Value = 21.27016493752229
$ws.Range("B3").Value = 22.24762660493843
$ws.Range("C3").Value = 8.339719073313539
$ws.Range("D3").Value = 12.64732452398467
$ws.Range("E3").Value = 12.23735975089206
$ws.Range("G3").Value = 58.76135607516459
$ws.Range("H3").Value = 21.73557944722167
$ws.Range("I3").Value = 33.13777545081675
$ws.Range("J3").Value = 7.669689668025995
$ws.Range("L3").Value = 13.2620300567123
$ws.Range("M3").Value = 20.33376262335443
$ws.Range("N3").Value = 21.342612008821
$ws.Range("B4").Value = 22.07718026780194
$ws.Range("C4").Value = 8.098226578982354
$ws.Range("D4").Value = 12.65942273304261
$ws.Range("E4").Value = 12.25032824361593
$ws.Range("G4").Value = 58.73695725333977
$ws.Range("H4").Value = 21.76375826862347
$ws.Range("I4").Value = 33.1922771324933
$ws.Range("J4").Value = 7.659665917269
$ws.Range("L4").Value = 13.2718720544607
$ws.Range("M4").Value = 20.30877501149429
$ws.Range("N4").Value = 21.38906456552273
$ws.Range("B5").Value = 22.00868254850037
$ws.Range("C5").Value = 7.998223914241187
$ws.Range("D5").Value = 12.66484135850401
$ws.Range("E5").Value = 12.25579095203504
$ws.Range("G5").Value = 58.7310784843006
$ws.Range("H5").Value = 21.77628118522683
$ws.Range("I5").Value = 33.2161332468224
$ws.Range("J5").Value = 7.65551663794426
$ws.Range("L5").Value = 13.27635760754669
$ws.Range("M5").Value = 20.29951252573457
$ws.Range("N5").Value = 21.40849131803534
$ws.Range("B6").Value = 21.99736890437514
$ws.Range("C6").Value = 7.981527355044127
$ws.Range("D6").Value = 12.66577061848772
$ws.Range("E6").Value = 12.25670879208251
$ws.Range("G6").Value = 58.73034770213303
$ws.Range("H6").Value = 21.77842334270655
$ws.Range("I6").Value = 33.22019383352971
$ws.Range("J6").Value = 7.654823756986633
$ws.Range("L6").Value = 13.27713111667017
$ws.Range("M6").Value = 20.29803026131186
$ws.Range("N6").Value = 21.4117471788975
$ws.Range("B7").Value = 22.07625248858537
$ws.Range("C7").Value = 8.096884133486466
$ws.Range("D7").Value = 12.65949383270768
$ws.Range("E7").Value = 12.25040119448542
$ws.Range("G7").Value = 58.73686151823787
$ws.Range("H7").Value = 21.76392294988706
$ws.Range("I7").Value = 33.19259220485753
$ws.Range("J7").Value = 7.659610219845119
$ws.Range("L7").Value = 13.27193062523029
$ws.Range("M7").Value = 20.30864636005388
$ws.Range("N7").Value = 21.3893245471919
$ws.Range("B8").Value = 22.43256407599166
$ws.Range("C8").Value = 8.59152125848493
$ws.Range("D8").Value = 12.6360846220828
$ws.Range("E8").Value = 12.22412804316213
$ws.Range("G8").Value = 58.80128448832427
$ws.Range("H8").Value = 21.70907038147379
$ws.Range("I8").Value = 33.08525661325718
$ws.Range("J8").Value = 7.68016534820559
$ws.Range("L8").Value = 13.2531501774785
$ws.Range("M8").Value = 20.36355273041567
$ws.Range("N8").Value = 21.29473685982601
$ws.Range("B9").Value = 23.15540363740435
$ws.Range("C9").Value = 9.500353608802035
$ws.Range("D9").Value = 12.60559645190195
$ws.Range("E9").Value = 12.17818420003544
$ws.Range("G9").Value = 59.05570013002591
$ws.Range("H9").Value = 21.63447154374701
$ws.Range("I9").Value = 32.92704379972538
$ws.Range("J9").Value = 7.718437549495623
$ws.Range("L9").Value = 13.23128368207391
$ws.Range("M9").Value = 20.4997031698923
$ws.Range("N9").Value = 21.12480325492409
$ws.Range("B10").Value = 23.69625396033085
$ws.Range("C10").Value = 10.12408692294239
$ws.Range("D10").Value = 12.59259403595051
$ws.Range("E10").Value = 12.14779909723481
$ws.Range("G10").Value = 59.32027603163028
$ws.Range("H10").Value = 21.59987457352116
$ws.Range("I10").Value = 32.84291953320297
$ws.Range("J10").Value = 7.745316592154632
$ws.Range("L10").Value = 13.22431851852878
$ws.Range("M10").Value = 20.61659954244581
$ws.Range("N10").Value = 21.00932489081068
$ws.Range("B11").Value = 23.94338168166124
$ws.Range("C11").Value = 10.39718627628087
$ws.Range("D11").Value = 12.58871746178839
$ws.Range("E11").Value = 12.13470139645309
$ws.Range("G11").Value = 59.45733383271736
$ws.Range("H11").Value = 21.58854428195269
$ws.Range("I11").Value = 32.81167123520603
$ws.Range("J11").Value = 7.757275823929974
$ws.Range("L11").Value = 13.22311969251083
$ws.Range("M11").Value = 20.67332085702228
$ws.Range("N11").Value = 20.95880235889913
$ws.Range("B12").Value = 24.03703372355661
$ws.Range("C12").Value = 10.49899794910398
$ws.Range("D12").Value = 12.58754225657682
$ws.Range("E12").Value = 12.12984535134765
$ws.Range("G12").Value = 59.51161536377969
$ws.Range("H12").Value = 21.58488890486268
$ws.Range("I12").Value = 32.8008513590545
$ws.Range("J12").Value = 7.761765887556322
$ws.Range("L12").Value = 13.22294824547741
$ws.Range("M12").Value = 20.69529812622817
$ws.Range("N12").Value = 20.93995796451476
$ws.Range("B13").Value = 24.01686232616388
$ws.Range("C13").Value = 10.47714344511159
$ws.Range("D13").Value = 12.58778234340507
$ws.Range("E13").Value = 12.13088657918097
$ws.Range("G13").Value = 59.49981935239247
$ws.Range("H13").Value = 21.58564788890178
$ws.Range("I13").Value = 32.80313649617874
$ws.Range("J13").Value = 7.760800591666528
$ws.Range("L13").Value = 13.22297261680958
$ws.Range("M13").Value = 20.69054296430485
$ws.Range("N13").Value = 20.94400368429953
$ws.Range("B14").Value = 23.95108557739353
$ws.Range("C14").Value = 10.4055949419202
$ws.Range("D14").Value = 12.58861491167276
$ws.Range("E14").Value = 12.13429980945302
$ws.Range("G14").Value = 59.46175202618527
$ws.Range("H14").Value = 21.58823081583574
$ws.Range("I14").Value = 32.81076075188921
$ws.Range("J14").Value = 7.757645998964584
$ws.Range("L14").Value = 13.22309992972644
$ws.Range("M14").Value = 20.67511902531418
$ws.Range("N14").Value = 20.9572462681538
$ws.Range("B15").Value = 23.91080200383359
$ws.Range("C15").Value = 10.36155835315475
$ws.Range("D15").Value = 12.58916299941594
$ws.Range("E15").Value = 12.13640401271967
$ws.Range("G15").Value = 59.43874404838677
$ws.Range("H15").Value = 21.58989568367447
$ws.Range("I15").Value = 32.81556287830839
$ws.Range("J15").Value = 7.755708685794106
$ws.Range("L15").Value = 13.22321468247019
$ws.Range("M15").Value = 20.6657359352467
$ws.Range("N15").Value = 20.96539511504635
$ws.Range("B16").Value = 23.68011809158971
$ws.Range("C16").Value = 10.10601844309929
$ws.Range("D16").Value = 12.59288836841078
$ws.Range("E16").Value = 12.14866960728893
$ws.Range("G16").Value = 59.31165334745438
$ws.Range("H16").Value = 21.60070385125822
$ws.Range("I16").Value = 32.84510326865542
$ws.Range("J16").Value = 7.744529625534077
$ws.Range("L16").Value = 13.22443644061022
$ws.Range("M16").Value = 20.61296302322447
$ws.Range("N16").Value = 21.01266695608034
$ws.Range("B17").Value = 23.53881835873434
$ws.Range("C17").Value = 9.946470574259667
$ws.Range("D17").Value = 12.59569560915287
$ws.Range("E17").Value = 12.15637944010263
$ws.Range("G17").Value = 59.23795114293081
$ws.Range("H17").Value = 21.60846430833005
$ws.Range("I17").Value = 32.86502609691325
$ws.Range("J17").Value = 7.737602895491786
$ws.Range("L17").Value = 13.22568995171587
$ws.Range("M17").Value = 20.58148838080381
$ws.Range("N17").Value = 21.04218020596257
$ws.Range("B18").Value = 23.45765440263494
$ws.Range("C18").Value = 9.853705050765608
$ws.Range("D18").Value = 12.59750213938785
$ws.Range("E18").Value = 12.16088216363661
$ws.Range("G18").Value = 59.1971330344356
$ws.Range("H18").Value = 21.61334276189455
$ws.Range("I18").Value = 32.87714575805132
$ws.Range("J18").Value = 7.733593585913121
$ws.Range("L18").Value = 13.22659638098477
$ws.Range("M18").Value = 20.56371936565212
$ws.Range("N18").Value = 21.05934464830826
$ws.Range("B19").Value = 23.43019489146185
$ws.Range("C19").Value = 9.822127347014943
$ws.Range("D19").Value = 12.59814676359214
$ws.Range("E19").Value = 12.16241844052309
$ws.Range("G19").Value = 59.18358354120834
$ws.Range("H19").Value = 21.61506572913996
$ws.Range("I19").Value = 32.88136260532758
$ws.Range("J19").Value = 7.732231770184369
$ws.Range("L19").Value = 13.22693515241594
$ws.Range("M19").Value = 20.55776085180274
$ws.Range("N19").Value = 21.06518877562821
$ws.Range("B20").Value = 23.55384942310601
$ws.Range("C20").Value = 9.963558567968869
$ws.Range("D20").Value = 12.59537691690713
$ws.Range("E20").Value = 12.15555165667109
$ws.Range("G20").Value = 59.24563416218158
$ws.Range("H20").Value = 21.60759524973695
$ws.Range("I20").Value = 32.86283688333715
$ws.Range("J20").Value = 7.738342870243697
$ws.Range("L20").Value = 13.22553732556993
$ws.Range("M20").Value = 20.58480438632
$ws.Range("N20").Value = 21.03901889963289
$ws.Range("B21").Value = 23.97040458332676
$ws.Range("C21").Value = 10.42665457532146
$ws.Range("D21").Value = 12.58836242365869
$ws.Range("E21").Value = 12.1332944481173
$ws.Range("G21").Value = 59.47286889485743
$ws.Range("H21").Value = 21.58745490048447
$ws.Range("I21").Value = 32.80849379797792
$ws.Range("J21").Value = 7.758573628539465
$ws.Range("L21").Value = 13.22305487343243
$ws.Range("M21").Value = 20.67963598614358
$ws.Range("N21").Value = 20.95334881261337
$ws.Range("B22").Value = 24.24301730969085
$ws.Range("C22").Value = 10.719931745709
$ws.Range("D22").Value = 12.58548430692316
$ws.Range("E22").Value = 12.11935272897109
$ws.Range("G22").Value = 59.63524435280549
$ws.Range("H22").Value = 21.57799463538044
$ws.Range("I22").Value = 32.77888452020805
$ws.Range("J22").Value = 7.771570350467324
$ws.Range("L22").Value = 13.22307880308232
$ws.Range("M22").Value = 20.74451163437924
$ws.Range("N22").Value = 20.89903289949104
$ws.Range("B23").Value = 24.09751315958397
$ws.Range("C23").Value = 10.5642846075424
$ws.Range("D23").Value = 12.58686442447676
$ws.Range("E23").Value = 12.12673850254412
$ws.Range("G23").Value = 59.54732074095384
$ws.Range("H23").Value = 21.58270460377722
$ws.Range("I23").Value = 32.79414591133401
$ws.Range("J23").Value = 7.764654390974989
$ws.Range("L23").Value = 13.22291565087792
$ws.Range("M23").Value = 20.70962512608033
$ws.Range("N23").Value = 20.92786961343229
$ws.Range("B24").Value = 23.54705365064507
$ws.Range("C24").Value = 9.955836317506014
$ws.Range("D24").Value = 12.59552039770227
$ws.Range("E24").Value = 12.1559256787506
$ws.Range("G24").Value = 59.24215582623763
$ws.Range("H24").Value = 21.60798685271212
$ws.Range("I24").Value = 32.86382455326486
$ws.Range("J24").Value = 7.738008411883615
$ws.Range("L24").Value = 13.22560574915309
$ws.Range("M24").Value = 20.5833042033502
$ws.Range("N24").Value = 21.040447512819
$ws.Range("B25").Value = 22.95778230682006
$ws.Range("C25").Value = 9.261793110064103
$ws.Range("D25").Value = 12.61219359833029
$ws.Range("E25").Value = 12.19001926095181
$ws.Range("G25").Value = 58.97319113595153
$ws.Range("H25").Value = 21.65111064794306
$ws.Range("I25").Value = 32.96422191281799
$ws.Range("J25").Value = 7.708301698820935
$ws.Range("L25").Value = 13.23559917916856
$ws.Range("M25").Value = 20.45987015025136
$ws.Range("N25").Value = 21.16912113461057
